$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row above the existing data row, shifting data to row 2
$ws.Rows.Item(1).Insert()

# Populate header row. Values are assigned in this specific order so that the
# generated shared-strings table indexes match the target (insertion order = index order):
# 0=a, 1=r, 2=d, 3=f, 4=sad, 5=doggy, 6=fight, 7=josh, 8=orange, 9=wario, 10=cool
$ws.Range("B1").Value = "a"
$ws.Range("A1").Value = "r"
$ws.Range("C1").Value = "d"
$ws.Range("D1").Value = "f"
$ws.Range("E1").Value = "sad"
$ws.Range("F1").Value = "doggy"
$ws.Range("G1").Value = "fight"
$ws.Range("H1").Value = "josh"
$ws.Range("I1").Value = "orange"
$ws.Range("J1").Value = "wario"
$ws.Range("K1").Value = "cool"

# Update the selection to match the target view state
$ws.Range("J6").Select()
